$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.597.27"
$ws.Range("E2").Value = "  -2.11%  "
$ws.Range("D3").Value = "2.648.34"
$ws.Range("E3").Value = "  -3.23%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "598.80"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.22%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "168.51"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.81%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("E8").Value = "  -0.61%  "
$ws.Range("D9").Value = "2.648.09"
$ws.Range("E9").Value = "  -3.20%  "
$ws.Range("E10").Value = "  -1.63%  "
$ws.Range("E11").Value = "  +1.94%  "
$ws.Range("E12").Value = "  -1.07%  "
$ws.Range("E13").Value = "  -2.22%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.13"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.52%  "
$ws.Range("D15").Value = "3.118.86"
$ws.Range("E15").Value = "  -3.96%  "
$ws.Range("E16").Value = "  -3.73%  "
$ws.Range("D17").Value = "67.398.63"
$ws.Range("E17").Value = "  -2.27%  "
$ws.Range("D18").Value = "2.642.31"
$ws.Range("E18").Value = "  -3.43%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.93"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.02%  "
$ws.Range("E20").Value = "  +2.59%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "364.10"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.27%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.42"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.22%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.81"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.07%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.03"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +8.33%  "
$ws.Range("E25").Value = "  -5.04%  "
$ws.Range("E26").Value = "  +0.07%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "70.78"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.34%  "
$ws.Range("E29").Value = "  -3.56%  "
$ws.Range("E30").Value = "  +0.03%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "557.54"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -5.65%  "
$ws.Range("E32").Value = "  -4.16%  "
$ws.Range("E33").Value = "  -4.10%  "
$ws.Range("E34").Value = "  -2.36%  "
$ws.Range("E35").Value = "  +0.04%  "
$ws.Range("E37").Value = "  -4.66%  "
$ws.Range("B38").Value = "EthereumClassic"
$ws.Range("C38").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "19.44"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.16%  "
$ws.Range("B39").Value = "Monero"
$ws.Range("C39").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "157.08"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.30%  "
$ws.Range("E40").Value = "  -2.70%  "
$ws.Range("E41").Value = "  -3.86%  "
$ws.Range("E43").Value = "  -0.43%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.54"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -5.63%  "
$ws.Range("E45").Value = "  +0.02%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "40.22"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.28%  "
$ws.Range("E47").Value = "  -3.36%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.598"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.52%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "154.03"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.50%  "
$ws.Range("E50").Value = "  -2.36%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.74"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.96%  "
